$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: extend header sequence with P1=14, Q1=15 (copy O1's format) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data rows 2..25 ---
for ($r = 2; $r -le 25; $r++) {
    # Swap I <-> K and M <-> O contents (I:1->2, K:2->1, M:1->2, O:2->1)
    $ws.Cells.Item($r, 9).Value  = 2   # I
    $ws.Cells.Item($r, 11).Value = 1   # K
    $ws.Cells.Item($r, 13).Value = 2   # M
    $ws.Cells.Item($r, 15).Value = 1   # O

    # New columns P and Q, value 2 each
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
